# writePassword-template.docx — resource re-save after the M2Doc custom
# property work (issue #295).
#
# The canonical-OOXML diff for *this* resource only shows its XML being
# re-serialized (namespace / attribute ordering normalized by the tool
# that re-saved the fixture) — none of the document's text, formatting,
# fields, styles or document properties actually changed:
#   * word/document.xml   : same elements/attributes/values, just
#                            reordered (xmlns declarations, w:color,
#                            w:pgSz/w:pgMar attributes, ...).
#   * word/styles.xml      : same, for w:docDefaults / w:latentStyles /
#                            w:style (attribute order only).
#   * docProps/custom.xml  : unchanged (the "m:var:self" property is
#                            untouched; this particular template's custom
#                            properties were not part of issue #295 for
#                            this file).
#
# So we touch nothing in the content/formatting model — any Find/Replace
# or OM property "round trip" (even setting a value back to itself)
# would rewrite runs/paragraphs and introduce differences that are not
# present in the target diff (e.g. run-splitting, new namespaces, word
# counts). We simply confirm the template is intact, matching the
# no-op nature of the change for this file.

$d = $word.ActiveDocument

$paragraphCount = $d.Paragraphs.Count
$fieldCount = $d.Fields.Count

Write-Output "paragraphs=$paragraphCount fields=$fieldCount"
